$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the header style from H1 onto the two new header cells, then set values
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)

$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

$ws.Range("I2").Value = 2
$ws.Range("J2").Value = 7

$ws.Range("I3").Value = 7
$ws.Range("J3").Value = 7

$ws.Range("I4").Value = 1
$ws.Range("J4").Value = 1
